# Weekly price-sheet update: a new (most recent) pricing pair is published
# for "Terminal La Palmera de La Serena - Acelga" and the whole historical
# block shifts down by one pair (2 rows), so the oldest pair that used to
# live at the bottom (rows 603:604) now survives as new rows 605:606.
#
# Implementation: insert two blank rows at the top of the block (469:470) -
# Excel shifts everything below down automatically, which reproduces the
# "cascade" for every other row exactly - then populate the two new rows
# with the new pricing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("469:470").Insert()

# Row 469 - "Primera" quality
$ws.Cells.Item(469, 1).Value = 8
$ws.Cells.Item(469, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(469, 3).Value = "Coquimbo"
$ws.Cells.Item(469, 4).Value = 44988
$ws.Cells.Item(469, 5).Value = 4
$ws.Cells.Item(469, 6).Value = 100112009
$ws.Cells.Item(469, 7).Value = "Acelga"
$ws.Cells.Item(469, 8).Value = "Sin especificar"
$ws.Cells.Item(469, 9).Value = "Primera"
$ws.Cells.Item(469, 10).Value = 2140
$ws.Cells.Item(469, 11).Value = 500
$ws.Cells.Item(469, 12).Value = 600
$ws.Cells.Item(469, 13).Value = 550
$ws.Cells.Item(469, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(469, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(469, 16).Value = 275
$ws.Cells.Item(469, 17).Value = 2
$ws.Cells.Item(469, 18).Value = "Hortaliza"

# Row 470 - "Segunda" quality
$ws.Cells.Item(470, 1).Value = 8
$ws.Cells.Item(470, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(470, 3).Value = "Coquimbo"
$ws.Cells.Item(470, 4).Value = 44988
$ws.Cells.Item(470, 5).Value = 4
$ws.Cells.Item(470, 6).Value = 100112009
$ws.Cells.Item(470, 7).Value = "Acelga"
$ws.Cells.Item(470, 8).Value = "Sin especificar"
$ws.Cells.Item(470, 9).Value = "Segunda"
$ws.Cells.Item(470, 10).Value = 1480
$ws.Cells.Item(470, 11).Value = 400
$ws.Cells.Item(470, 12).Value = 450
$ws.Cells.Item(470, 13).Value = 425
$ws.Cells.Item(470, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(470, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(470, 16).Value = 212
$ws.Cells.Item(470, 17).Value = 2
$ws.Cells.Item(470, 18).Value = "Hortaliza"
